# Updated cryptos list on Thu May  2 02:54:42 UTC 2024 with GitHub Actions
#
# Sets cell values as plain text (matching the source workbook's inlineStr
# cells) without letting Excel auto-coerce numeric-looking strings (e.g.
# "547.88") into real numbers, and without leaving a stray NumberFormat
# behind on the cell once we're done.
function Set-TextValue {
    param($ws, $addr, $val)
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
Set-TextValue $ws "D2" "57.577.47"
Set-TextValue $ws "E2" "  -3.93%  "

# Row 3 - Ethereum
Set-TextValue $ws "D3" "2.924.26"
Set-TextValue $ws "E3" "  -2.17%  "

# Row 4 - TetherUSD
Set-TextValue $ws "E4" "  -0.06%  "

# Row 5 - BNB
Set-TextValue $ws "D5" "547.88"
Set-TextValue $ws "E5" "  -3.93%  "

# Row 6 - Solana
Set-TextValue $ws "D6" "129.98"
Set-TextValue $ws "E6" "  +3.74%  "

# Row 7
Set-TextValue $ws "E7" "  -0.08%  "

# Row 8
Set-TextValue $ws "D8" "0.510"
Set-TextValue $ws "E8" "  +1.78%  "

# Row 9
Set-TextValue $ws "D9" "2.924.36"
Set-TextValue $ws "E9" "  -2.05%  "

# Row 10
Set-TextValue $ws "E10" "  -3.18%  "

# Row 11
Set-TextValue $ws "D11" "4.76"
Set-TextValue $ws "E11" "  -5.82%  "

# Row 12
Set-TextValue $ws "D12" "0.446"
Set-TextValue $ws "E12" "  +1.37%  "

# Row 13
Set-TextValue $ws "D13" "0.0000220"
Set-TextValue $ws "E13" "  +0.27%  "

# Row 14
Set-TextValue $ws "D14" "32.82"
Set-TextValue $ws "E14" "  +0.95%  "

# Row 15
Set-TextValue $ws "E15" "  -0.05%  "

# Row 16
Set-TextValue $ws "D16" "3.400.23"
Set-TextValue $ws "E16" "  -2.31%  "

# Row 17
Set-TextValue $ws "D17" "6.83"
Set-TextValue $ws "E17" "  +5.83%  "

# Row 18
Set-TextValue $ws "D18" "2.920.59"
Set-TextValue $ws "E18" "  -2.38%  "

# Row 19
Set-TextValue $ws "D19" "57.542.96"
Set-TextValue $ws "E19" "  -4.03%  "

# Row 20
Set-TextValue $ws "D20" "417.04"
Set-TextValue $ws "E20" "  -2.08%  "

# Row 21
Set-TextValue $ws "D21" "13.12"
Set-TextValue $ws "E21" "  +0.27%  "

# Row 22
Set-TextValue $ws "D22" "0.682"
Set-TextValue $ws "E22" "  +2.17%  "

# Row 23
Set-TextValue $ws "D23" "6.95"
Set-TextValue $ws "E23" "  -1.14%  "

# Row 24
Set-TextValue $ws "D24" "12.98"
Set-TextValue $ws "E24" "  +0.44%  "

# Row 25
Set-TextValue $ws "D25" "79.74"
Set-TextValue $ws "E25" "  +0.71%  "

# Row 27
Set-TextValue $ws "E27" "  -0.08%  "

# Row 28
Set-TextValue $ws "D28" "2.46"
Set-TextValue $ws "E28" "  -2.46%  "

# Row 29
Set-TextValue $ws "D29" "7.43"
Set-TextValue $ws "E29" "  +2.31%  "

# Row 30
Set-TextValue $ws "D30" "1.98"
Set-TextValue $ws "E30" "  +1.77%  "

# Row 31
Set-TextValue $ws "D31" "25.16"
Set-TextValue $ws "E31" "  +0.07%  "

# Row 32
Set-TextValue $ws "D32" "5.95"
Set-TextValue $ws "E32" "  -3.34%  "

# Row 33
Set-TextValue $ws "D33" "0.0970"
Set-TextValue $ws "E33" "  +3.51%  "

# Row 34
Set-TextValue $ws "D34" "5.63"
Set-TextValue $ws "E34" "  +0.44%  "

# Row 35
Set-TextValue $ws "D35" "0.938"
Set-TextValue $ws "E35" "  +1.13%  "

# Row 36
Set-TextValue $ws "D36" "2.06"
Set-TextValue $ws "E36" "  +0.56%  "

# Row 37
Set-TextValue $ws "D37" "47.85"
Set-TextValue $ws "E37" "  -4.33%  "

# Row 38
Set-TextValue $ws "D38" "8.69"
Set-TextValue $ws "E38" "  +4.08%  "

# Row 39
Set-TextValue $ws "D39" "0.0₃0676"
Set-TextValue $ws "E39" "  +2.49%  "

# Row 40
Set-TextValue $ws "D40" "2.55"
Set-TextValue $ws "E40" "  +3.97%  "

# Row 41 - was Bittensor, now Kaspa
Set-TextValue $ws "B41" "Kaspa"
Set-TextValue $ws "C41" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws "D41" "0.107"
Set-TextValue $ws "E41" "  -1.38%  "

# Row 42 - was Kaspa, now Bittensor
Set-TextValue $ws "B42" "Bittensor"
Set-TextValue $ws "C42" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws "D42" "374.57"
Set-TextValue $ws "E42" "  -0.27%  "

# Row 43
Set-TextValue $ws "D43" "0.0344"
Set-TextValue $ws "E43" "  -2.73%  "

# Row 44
Set-TextValue $ws "D44" "2.658.68"
Set-TextValue $ws "E44" "  -0.32%  "

# Row 46 - was Monero, now TheGraph
Set-TextValue $ws "B46" "TheGraph"
Set-TextValue $ws "C46" "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws "D46" "0.238"
Set-TextValue $ws "E46" "  +1.49%  "

# Row 47 - was TheGraph, now Monero
Set-TextValue $ws "B47" "Monero"
Set-TextValue $ws "C47" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws "D47" "122.08"
Set-TextValue $ws "E47" "  +1.62%  "

# Row 48
Set-TextValue $ws "D48" "0.109"
Set-TextValue $ws "E48" "  +1.62%  "

# Row 49
Set-TextValue $ws "D49" "1.96"
Set-TextValue $ws "E49" "  -1.52%  "

# Row 50
Set-TextValue $ws "D50" "23.11"
Set-TextValue $ws "E50" "  -1.97%  "

# Row 51
Set-TextValue $ws "D51" "2.00"
Set-TextValue $ws "E51" "  -0.20%  "
